# Add the "CreateNewContact" worksheet (new test case 3) to the freeCRM
# workbook, fill in its header/data row, wire up the two mailto hyperlinks,
# and make it the active sheet - matching the existing
# Authentication / CreateNewCompany sheets' look & feel.

$wb = $excel.ActiveWorkbook

# --- Fix up CreateNewCompany's lingering selection before we move away from it ---
$ws2 = $wb.Worksheets.Item("CreateNewCompany")
$ws2.Activate()
$ws2.Range("A2:B2").Select()

# --- Create the new sheet after the last existing sheet ---
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws3 = $wb.Worksheets.Add($null, $lastSheet)
$ws3.Name = "CreateNewContact"

# --- Column widths (approximate character widths matching the source form) ---
$ws3.Columns.Item(1).ColumnWidth = 25.666666666666664
$ws3.Columns.Item(2).ColumnWidth = 16.333333333333336
$ws3.Columns.Item(6).ColumnWidth = 10.0
$ws3.Columns.Item(7).ColumnWidth = 25.833333333333336
$ws3.Columns.Item(8).ColumnWidth = 10.0
$ws3.Columns.Item(9).ColumnWidth = 11.833333333333332
$ws3.Columns.Item(10).ColumnWidth = 25.0
$ws3.Columns.Item(11).ColumnWidth = 26.166666666666664
$ws3.Columns.Item(12).ColumnWidth = 12.333333333333332
$ws3.Columns.Item(13).ColumnWidth = 11.333333333333332
$ws3.Columns.Item(24).ColumnWidth = 11.666666666666668
$ws3.Columns.Item(25).ColumnWidth = 13.166666666666668

# --- Pre-format the data row as Text so every value (including numeric-looking
#     ones like IDs/years/zip codes) lands as a shared string, not a number ---
$ws3.Range("A2:X2").NumberFormat = "@"

# --- Header / first data row values ---
# (written in the same order the strings were first entered in the source
# workbook, so new shared-string indices line up with the canonical file)
$ws3.Range("A2").Value = "adilkhaleque429@gmail.com"
$ws3.Range("B2").Value = "Testunbound6F"
$ws3.Range("C2").Value = "Simon"
$ws3.Range("D2").Value = "Winter"
$ws3.Range("E2").Value = "Company1"
$ws3.Range("F2").Value = "adilkhaleque429@gmail.com"
$ws3.Range("L2").Value = "11 Hydrant St"
$ws3.Range("M2").Value = "Philadelphia"
$ws3.Range("N2").Value = "PA"
$ws3.Range("P2").Value = "2111111111"
$ws3.Range("Q2").Value = "Mobile"
$ws3.Range("R2").Value = "Manager"
$ws3.Range("S2").Value = "IT"
$ws3.Range("H2").Value = "1"
$ws3.Range("T2").Value = "2"
$ws3.Range("I2").Value = "1"
$ws3.Range("W2").Value = "1992"
$ws3.Range("G2").Value = "Business"
$ws3.Range("O2").Value = "12121"
$ws3.Range("X2").Value = "Simon Winter"
$ws3.Range("U2").Value = "1"
$ws3.Range("V2").Value = "1"
$ws3.Range("K2").Value = "152"
$ws3.Range("J2").Value = "New contact for Company1."

# --- Hyperlinks for the two email-address cells (also promotes them to the
#     combined Hyperlink+Text style used by the source workbook) ---
$ws3.Hyperlinks.Add($ws3.Range("A2"), "mailto:adilkhaleque429@gmail.com")
$ws3.Hyperlinks.Add($ws3.Range("F2"), "mailto:adilkhaleque429@gmail.com")

# --- Selection / activation ---
$ws3.Range("J2").Select()
$ws3.Activate()
